# Fix issue with unknown cell type
# - Adds a "float" column (H) and a "str" column (I, formula-based) between
#   the existing "int" column and the "Empty column" column (moved to J/K).
# - Swaps the two email addresses (and their hyperlinks) between John and Jean.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header row (row 1)
# ---------------------------------------------------------------------------
# H1 used to be the (now relocated) "Empty column" header; it becomes the new
# "Column with float" header and picks up the same font used by the
# hyperlink cells (blue Droid Sans 10pt already present in the workbook).
$ws.Range("H1").Value = "Column with float"
$ws.Range("H1").Font.Name = "Droid Sans"
$ws.Range("H1").Font.Size = 10
$ws.Range("H1").Font.Color = 16711680

$ws.Range("I1").Value = "Column with str"
$ws.Range("J1").Value = "Empty column"

# ---------------------------------------------------------------------------
# Row 2 (Jean / Lefebvre record)
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = "john.lefebvre@mail.com"

# New "float" column: stored as text (so it round-trips exactly), not as a
# double, to avoid floating point artefacts such as 5073.0100000000002.
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "5073.01"
$ws.Range("H2").NumberFormat = "General"

# The old "str" column content moves over to K2 ...
$ws.Range("K2").Value = "dsdsqd"
# ... and I2 becomes a formula reading back the float column as text.
$ws.Range("I2").Formula = "=LEFT(H2,10)"

$ws.Rows.Item(2).RowHeight = 13.8

# ---------------------------------------------------------------------------
# Row 4 (John / Smith record)
# ---------------------------------------------------------------------------
$ws.Range("C4").Value = "jean.smith@mail.com"

$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "0001.01"
$ws.Range("H4").NumberFormat = "General"

$ws.Range("K4").Value = "qdqsdq"
$ws.Range("I4").Formula = "=LEFT(H4,10)"

$ws.Rows.Item(4).RowHeight = 13.8

# Give the two new float cells their own font (Calibri 11, black) - done as a
# single multi-cell assignment so both cells end up sharing one style.
$floatRange = $ws.Range("H2,H4")
$floatRange.Font.Name = "Calibri"
$floatRange.Font.Size = 11
$floatRange.Font.Color = 0

# ---------------------------------------------------------------------------
# Hyperlinks: re-target rather than append, so John's/Jean's mail links now
# point at the swapped addresses.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:john.lefebvre@mail.com", "", "", "john.lefebvre@mail.com")
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:jean.smith@mail.com", "", "", "jean.smith@mail.com")

# Matches the author's final selection position in the fixture.
$ws.Range("I5").Select()
